# BOT; UPDATE DATA (#1654)
# Update the PCR/infection-tracking figures on the "all" and "kobe"
# sheets, then restore the UI selection state (active sheet/tab,
# selected cells) to match the authored workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "all" sheet - cumulative recovered (column H) revisions, plus
#    one cumulative tested (column C) revision on row 40.
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Range("H26").Value = 165
$wsAll.Range("H27").Value = 165
$wsAll.Range("H28").Value = 176
$wsAll.Range("H29").Value = 179
$wsAll.Range("H30").Value = 184
$wsAll.Range("H31").Value = 189
$wsAll.Range("H32").Value = 196
$wsAll.Range("H33").Value = 198
$wsAll.Range("H34").Value = 202
$wsAll.Range("H35").Value = 213
$wsAll.Range("H36").Value = 213
$wsAll.Range("H37").Value = 219
$wsAll.Range("H38").Value = 223
$wsAll.Range("H39").Value = 227
$wsAll.Range("C40").Value = 281
$wsAll.Range("H40").Value = 228
$wsAll.Range("H41").Value = 229

# ---------------------------------------------------------------
# 2. "kobe" sheet - cumulative tested (column E) and cumulative
#    recovered (column J) revisions.
# ---------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("J81").Value = 158
$wsKobe.Range("E82").Value = 268
$wsKobe.Range("J82").Value = 158
$wsKobe.Range("J83").Value = 169
$wsKobe.Range("J84").Value = 172
$wsKobe.Range("J85").Value = 177
$wsKobe.Range("J86").Value = 180
$wsKobe.Range("J87").Value = 187
$wsKobe.Range("J88").Value = 189
$wsKobe.Range("J89").Value = 193
$wsKobe.Range("J90").Value = 204
$wsKobe.Range("J91").Value = 204
$wsKobe.Range("J92").Value = 210
$wsKobe.Range("J93").Value = 214
$wsKobe.Range("J94").Value = 218
$wsKobe.Range("J95").Value = 219
$wsKobe.Range("J96").Value = 220

# ---------------------------------------------------------------
# 3. Restore each sheet's own selection, then finish with "all"
#    active/selected so the workbook re-opens on that tab (matches
#    the dropped workbookView activeTab="2" -> defaults to tab 0).
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Range("A72").Select()

$wsKobe.Activate()
$wsKobe.Range("A80").Select()

$wsAll.Activate()
$wsAll.Range("L35").Select()
